# Apply the "new timing" update (isaqb-org/curriculum-foundation issue #99).
#
# The raw timing figures live on the "Times" sheet (columns B/C, with column D
# a per-row total formula and D8 a grand total). The two graphic sheets
# ("EN-Graphic" / "DE-Grafik") just pull column D of "Times" via formulas
# (=Times!D3, =Times!D4, ...), and the pie charts on those sheets are driven
# from that same graphic-sheet data — so editing the three "Times" source
# cells below is enough for every dependent formula (and, inside a real
# Excel session, the chart caches) to follow along.

$wb = $excel.ActiveWorkbook

# --- 1. Update the raw timings on the "Times" sheet -----------------------
$wsTimes = $wb.Worksheets.Item("Times")
$wsTimes.Range("B3").Value = 300   # was 330  -> Times!D3 (=B3+C3) becomes 420
$wsTimes.Range("B4").Value = 150   # was 180  -> Times!D4 (=B4+C4) becomes 240
$wsTimes.Range("C4").Value = 90    # was 120  -> Times!D4 (=B4+C4) becomes 240

# --- 2. Match the saved selection/active-sheet state -----------------------
# "Times" sheet selection moves from A1:D7 to the single cell C4.
$wsTimes.Activate()
$wsTimes.Range("C4").Select()

# The workbook's active tab moves from "DE-Grafik" back to "EN-Graphic".
$wsEN = $wb.Worksheets.Item("EN-Graphic")
$wsEN.Activate()
